$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: harrystyles / 123456 (text) / Istutor=False / Score=0
$ws.Range("A13").Value = "harrystyles"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "123456"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = 0

# Row 14: eltonjohn / 123456 (text) / Istutor=True / Score=0
$ws.Range("A14").Value = "eltonjohn"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "123456"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = $true
$ws.Range("D14").Value = 0
